$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Export" sheet is a flat dump of account balances (Conta / Nome / Saldo)
# with one row per account, sorted descending by balance. This edit:
#   - removes the RACHEL (004940699) row entirely
#   - doubles OLGA's (008004799) balance from 60000 to 120000
#   - removes the PAULA / JULIA / BERNARDO / REJANE rows entirely
#   - re-sorts the HELOISA..CINCO block above CARLOS (keeping their own values)
#   - drops CARLOS's (004386464) balance from 27001.47 to 8001.47, which moves
#     his row down below the HELOISA..CINCO block to keep the sheet sorted
#
# Net effect on rows 5-20 (1-indexed, header = row 1): 5 rows disappear
# (RACHEL, PAULA, JULIA, BERNARDO, REJANE) and the remaining rows 4-15 end up
# holding a specific new sequence of (Conta, Nome, Saldo) triples.
#
# Delete the five now-absent rows first, bottom-to-top so row numbers of the
# rows still to be deleted don't shift.
$ws.Rows(20).Delete()   # REJANE     004550605
$ws.Rows(19).Delete()   # BERNARDO   005262440
$ws.Rows(18).Delete()   # JULIA      004265173
$ws.Rows(12).Delete()   # PAULA      005920340
$ws.Rows(5).Delete()    # RACHEL     004940699

# After those deletions, rows 4-15 contain (in this order):
#   4  OLGA      60000
#   5  ROGERIO   55152.2
#   6  GILSON    38195.95
#   7  CLAUDIA   35028.02
#   8  KELLY     31938.06
#   9  CARLOS    27001.47
#   10 HELOISA   24205.55
#   11 JOSE      23156.83
#   12 RENATO    23091.26
#   13 EDUARDO   22080.72
#   14 BRUNO     15414.17
#   15 CINCO     14455.12
#
# Rewrite that contiguous block with the final values/order: OLGA's new
# balance, and the HELOISA..CINCO block moved above CARLOS (whose balance
# also changes).
$data = @(
    @("008004799", "OLGA",    120000),
    @("004487016", "ROGERIO", 55152.2),
    @("004474776", "GILSON",  38195.95),
    @("002697806", "CLAUDIA", 35028.02),
    @("004556974", "KELLY",   31938.06),
    @("004468717", "HELOISA", 24205.55),
    @("004453157", "JOSE",    23156.83),
    @("004862672", "RENATO",  23091.26),
    @("004461070", "EDUARDO", 22080.72),
    @("004515341", "BRUNO",   15414.17),
    @("004581652", "CINCO",   14455.12),
    @("004386464", "CARLOS",  8001.47)
)

$startRow = 4
$endRow = $startRow + $data.Count - 1

# Account numbers ("Conta") are zero-padded numeric strings (e.g. "008004799").
# Force column A to Text format first so Excel doesn't reinterpret the
# assigned strings as numbers and strip the leading zeros.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
